# GoodInfo_v2 - 2021.12.08 完成
# Append 8 new ticker rows (95-102) to the tracking sheet, matching the
# style of the existing data rows, and move the active selection to the
# new last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New (code, name) pairs to append below the existing last row (94).
$newRows = @(
    @{Code = 3624; Name = "光頡"},
    @{Code = 2442; Name = "新美齊"},
    @{Code = 5356; Name = "協益"},
    @{Code = 2006; Name = "東和鋼鐵"},
    @{Code = 2303; Name = "聯電"},
    @{Code = 2405; Name = "浩鑫"},
    @{Code = 3294; Name = "英濟"},
    @{Code = 8104; Name = "錸寶"}
)

# Find the current last used row in column A (xlUp from the bottom of the
# sheet), then grab its formatting so the new rows keep the same style as
# the rest of the table.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$templateStyle = $ws.Range("A" + $lastRow + ":B" + $lastRow).Style

foreach ($row in $newRows) {
    $lastRow = $lastRow + 1

    $rowRange = $ws.Range("A" + $lastRow + ":B" + $lastRow)
    $rowRange.Style = $templateStyle

    $ws.Cells.Item($lastRow, 1).Value = $row.Code
    $ws.Cells.Item($lastRow, 2).Value = $row.Name
}

# Match the saved selection/active cell from the source workbook.
[void]$ws.Range("A" + $lastRow + ":B" + $lastRow).Select()
